# Apply the motilal_portfolio_change_engine update:
# Insert a new "Industry" column between "Stock Name" (B) and "Mutual Fund" (C),
# shifting the existing C:I columns to D:J, and populate the new column
# with industry classifications for each holding row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C; this shifts existing C:I -> D:J and copies formatting.
$ws.Range("C1").EntireColumn.Insert()

# Set the new header cell, matching the style/formatting of the other header cells
# (bold/center/border), by copying formats from the adjacent header cell.
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C1").Value = "Industry"

# Map of row -> Industry classification.
$industries = @{
    2  = "IT - Software"
    3  = "Telecom - Services"
    4  = "IT - Software"
    5  = "Capital Markets"
    6  = "IT - Software"
    7  = "Retailing"
    8  = "Retailing"
    9  = "Financial Technology (Fintech)"
    10 = "IT - Software"
    11 = "IT - Software"
    12 = "Retailing"
    13 = "Retailing"
    14 = "Retailing"
    15 = "Commercial Services & Supplies"
    16 = "Capital Markets"
    17 = "Electrical Equipment"
    18 = "IT - Services"
    19 = "Healthcare Services"
    20 = "IT - Software"
    21 = "Consumer Durables"
    22 = "IT - Services"
    23 = "Other Consumer Services"
    24 = "IT - Software"
    25 = "IT - Software"
    26 = "IT - Services"
    27 = "Electrical Equipment"
    28 = "Retailing"
    29 = "Capital Markets"
    30 = "Industrial Manufacturing"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
